$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 6.201571333333333
$ws.Cells.Item(2, 8).Value = 18.604714
$ws.Cells.Item(2, 9).Value = 0.05221490529364391
$ws.Cells.Item(2, 10).Value = 0.07406232529850043
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 18.59297266666667
$ws.Cells.Item(2, 14).Value = 55.778918
$ws.Cells.Item(2, 15).Value = 0.1361024717868078
$ws.Cells.Item(2, 16).Value = 0.1378163998421381
$ws.Cells.Item(2, 17).Value = 115.3056462910502
$ws.Cells.Item(2, 18).Value = 1037.750816619452
$ws.Cells.Item(2, 19).Value = 0.007106577674579014
$ws.Cells.Item(2, 20).Value = 0.01020700303657664

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 6.201571333333333
$ws.Cells.Item(3, 8).Value = 18.604714
$ws.Cells.Item(3, 9).Value = 0.05221490529364391
$ws.Cells.Item(3, 10).Value = 0.07406232529850043
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 6.746562
$ws.Cells.Item(3, 14).Value = 20.239686
$ws.Cells.Item(3, 15).Value = 0.04938552757134602
$ws.Cells.Item(3, 16).Value = 0.05000743575655815
$ws.Cells.Item(3, 17).Value = 41.83928549775599
$ws.Cells.Item(3, 18).Value = 376.5535694798039
$ws.Cells.Item(3, 19).Value = 0.002578660645014473
$ws.Cells.Item(3, 20).Value = 0.003703666974346071

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 6.201571333333333
$ws.Cells.Item(4, 8).Value = 18.604714
$ws.Cells.Item(4, 9).Value = 0.05221490529364391
$ws.Cells.Item(4, 10).Value = 0.07406232529850043
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 43.343503
$ws.Cells.Item(4, 14).Value = 130.030509
$ws.Cells.Item(4, 15).Value = 0.3172788988591848
$ws.Cells.Item(4, 16).Value = 0.3212743678538321
$ws.Cells.Item(4, 17).Value = 268.7978256910473
$ws.Cells.Item(4, 18).Value = 2419.180431219426
$ws.Cells.Item(4, 19).Value = 0.01656668765560396
$ws.Cells.Item(4, 20).Value = 0.0237943267420606

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 6.201571333333333
$ws.Cells.Item(5, 8).Value = 18.604714
$ws.Cells.Item(5, 9).Value = 0.05221490529364391
$ws.Cells.Item(5, 10).Value = 0.07406232529850043
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 62.83028666666667
$ws.Cells.Item(5, 14).Value = 188.49086
$ws.Cells.Item(5, 15).Value = 0.4599241590742428
$ws.Cells.Item(5, 16).Value = 0.4657159489602949
$ws.Cells.Item(5, 17).Value = 389.6465046571155
$ws.Cells.Item(5, 18).Value = 3506.81854191404
$ws.Cells.Item(5, 19).Value = 0.0240148964083204
$ws.Cells.Item(5, 20).Value = 0.03449200610859719

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 6.201571333333333
$ws.Cells.Item(6, 8).Value = 18.604714
$ws.Cells.Item(6, 9).Value = 0.05221490529364391
$ws.Cells.Item(6, 10).Value = 0.07406232529850043
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 5.0967785
$ws.Cells.Item(6, 14).Value = 10.193557
$ws.Cells.Item(6, 15).Value = 0.03730894270841855
$ws.Cells.Item(6, 16).Value = 0.02518584758717668
$ws.Cells.Item(6, 17).Value = 31.60803543794966
$ws.Cells.Item(6, 18).Value = 189.648212627698
$ws.Cells.Item(6, 19).Value = 0.001948082910126061
$ws.Cells.Item(6, 20).Value = 0.001865322436919931

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 7.461641333333333
$ws.Cells.Item(7, 8).Value = 22.384924
$ws.Cells.Item(7, 9).Value = 0.06282422221945559
$ws.Cells.Item(7, 10).Value = 0.0891107233935555
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.59297266666667
$ws.Cells.Item(7, 14).Value = 55.778918
$ws.Cells.Item(7, 15).Value = 0.1361024717868078
$ws.Cells.Item(7, 16).Value = 0.1378163998421381
$ws.Cells.Item(7, 17).Value = 138.7340933591369
$ws.Cells.Item(7, 18).Value = 1248.606840232232
$ws.Cells.Item(7, 19).Value = 0.008550531932151601
$ws.Cells.Item(7, 20).Value = 0.01228091908542842

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 7.461641333333333
$ws.Cells.Item(8, 8).Value = 22.384924
$ws.Cells.Item(8, 9).Value = 0.06282422221945559
$ws.Cells.Item(8, 10).Value = 0.0891107233935555
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 6.746562
$ws.Cells.Item(8, 14).Value = 20.239686
$ws.Cells.Item(8, 15).Value = 0.04938552757134602
$ws.Cells.Item(8, 16).Value = 0.05000743575655815
$ws.Cells.Item(8, 17).Value = 50.340425877096
$ws.Cells.Item(8, 18).Value = 453.0638328938639
$ws.Cells.Item(8, 19).Value = 0.003102607358567294
$ws.Cells.Item(8, 20).Value = 0.00445619877532365

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 7.461641333333333
$ws.Cells.Item(9, 8).Value = 22.384924
$ws.Cells.Item(9, 9).Value = 0.06282422221945559
$ws.Cells.Item(9, 10).Value = 0.0891107233935555
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 43.343503
$ws.Cells.Item(9, 14).Value = 130.030509
$ws.Cells.Item(9, 15).Value = 0.3172788988591848
$ws.Cells.Item(9, 16).Value = 0.3212743678538321
$ws.Cells.Item(9, 17).Value = 323.4136735162573
$ws.Cells.Item(9, 18).Value = 2910.723061646316
$ws.Cells.Item(9, 19).Value = 0.0199328000474736
$ws.Cells.Item(9, 20).Value = 0.02862899132726223

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 7.461641333333333
$ws.Cells.Item(10, 8).Value = 22.384924
$ws.Cells.Item(10, 9).Value = 0.06282422221945559
$ws.Cells.Item(10, 10).Value = 0.0891107233935555
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 62.83028666666667
$ws.Cells.Item(10, 14).Value = 188.49086
$ws.Cells.Item(10, 15).Value = 0.4599241590742428
$ws.Cells.Item(10, 16).Value = 0.4657159489602949
$ws.Cells.Item(10, 17).Value = 468.8170639771822
$ws.Cells.Item(10, 18).Value = 4219.353575794639
$ws.Cells.Item(10, 19).Value = 0.02889437757377647
$ws.Cells.Item(10, 20).Value = 0.04150028510776805

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 7.461641333333333
$ws.Cells.Item(11, 8).Value = 22.384924
$ws.Cells.Item(11, 9).Value = 0.06282422221945559
$ws.Cells.Item(11, 10).Value = 0.0891107233935555
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 5.0967785
$ws.Cells.Item(11, 14).Value = 10.193557
$ws.Cells.Item(11, 15).Value = 0.03730894270841855
$ws.Cells.Item(11, 16).Value = 0.02518584758717668
$ws.Cells.Item(11, 17).Value = 38.03033312244467
$ws.Cells.Item(11, 18).Value = 228.181998734668
$ws.Cells.Item(11, 19).Value = 0.002343905307486625
$ws.Cells.Item(11, 20).Value = 0.002244329097773148

# Row 12
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 7).Value = 105.106922
$ws.Cells.Item(12, 8).Value = 210.213844
$ws.Cells.Item(12, 9).Value = 0.8849608724869005
$ws.Cells.Item(12, 10).Value = 0.836826951307944
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 18.59297266666667
$ws.Cells.Item(12, 14).Value = 55.778918
$ws.Cells.Item(12, 15).Value = 0.1361024717868078
$ws.Cells.Item(12, 16).Value = 0.1378163998421381
$ws.Cells.Item(12, 17).Value = 1954.250127823465
$ws.Cells.Item(12, 18).Value = 11725.50076694079
$ws.Cells.Item(12, 19).Value = 0.1204453621800772
$ws.Cells.Item(12, 20).Value = 0.1153284777201331

# Row 13
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 7).Value = 105.106922
$ws.Cells.Item(13, 8).Value = 210.213844
$ws.Cells.Item(13, 9).Value = 0.8849608724869005
$ws.Cells.Item(13, 10).Value = 0.836826951307944
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 6.746562
$ws.Cells.Item(13, 14).Value = 20.239686
$ws.Cells.Item(13, 15).Value = 0.04938552757134602
$ws.Cells.Item(13, 16).Value = 0.05000743575655815
$ws.Cells.Item(13, 17).Value = 709.1103659021639
$ws.Cells.Item(13, 18).Value = 4254.662195412983
$ws.Cells.Item(13, 19).Value = 0.04370425956776425
$ws.Cells.Item(13, 20).Value = 0.04184757000688842

# Row 14
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 7).Value = 105.106922
$ws.Cells.Item(14, 8).Value = 210.213844
$ws.Cells.Item(14, 9).Value = 0.8849608724869005
$ws.Cells.Item(14, 10).Value = 0.836826951307944
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 43.343503
$ws.Cells.Item(14, 14).Value = 130.030509
$ws.Cells.Item(14, 15).Value = 0.3172788988591848
$ws.Cells.Item(14, 16).Value = 0.3212743678538321
$ws.Cells.Item(14, 17).Value = 4555.702189027766
$ws.Cells.Item(14, 18).Value = 27334.2131341666
$ws.Cells.Item(14, 19).Value = 0.2807794111561072
$ws.Cells.Item(14, 20).Value = 0.2688510497845092

# Row 15
$ws.Cells.Item(15, 5).Value = 2
$ws.Cells.Item(15, 7).Value = 105.106922
$ws.Cells.Item(15, 8).Value = 210.213844
$ws.Cells.Item(15, 9).Value = 0.8849608724869005
$ws.Cells.Item(15, 10).Value = 0.836826951307944
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 62.83028666666667
$ws.Cells.Item(15, 14).Value = 188.49086
$ws.Cells.Item(15, 15).Value = 0.4599241590742428
$ws.Cells.Item(15, 16).Value = 0.4657159489602949
$ws.Cells.Item(15, 17).Value = 6603.898039910973
$ws.Cells.Item(15, 18).Value = 39623.38823946584
$ws.Cells.Item(15, 19).Value = 0.4070148850921459
$ws.Cells.Item(15, 20).Value = 0.3897236577439296

# Row 16
$ws.Cells.Item(16, 5).Value = 2
$ws.Cells.Item(16, 7).Value = 105.106922
$ws.Cells.Item(16, 8).Value = 210.213844
$ws.Cells.Item(16, 9).Value = 0.8849608724869005
$ws.Cells.Item(16, 10).Value = 0.836826951307944
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 5.0967785
$ws.Cells.Item(16, 14).Value = 10.193557
$ws.Cells.Item(16, 15).Value = 0.03730894270841855
$ws.Cells.Item(16, 16).Value = 0.02518584758717668
$ws.Cells.Item(16, 17).Value = 535.706700250777
$ws.Cells.Item(16, 18).Value = 2142.826801003108
$ws.Cells.Item(16, 19).Value = 0.03301695449080586
$ws.Cells.Item(16, 20).Value = 0.0210761960524836
